$d = $word.ActiveDocument

# The document mentions the "Fakturoid" add-on once early on (right after
# "doplnek"/"doplněk"); a second, unrelated mention later refers to the
# company "Fakturoid s.r.o.". We only want to touch the first one, turning
# it into "Fakturoid by Kulhánek" (keeping the existing 10pt formatting),
# so search forward from the very start of the document and replace only
# that single occurrence.
$d.Content.Find.Execute("Fakturoid", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Fakturoid by Kulhánek", 1) | Out-Null
